$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Comment text (T12): clarify this is "before standardizing" ---
$ws.Range("T12").Value = "Currently, what is coming out is like this in the mixed pattern, before standardizing:"

# --- New content appended below the existing table (rows 14-18) ---
$ws.Range("W14").Value = "This looks good!"

$ws.Range("T16").Value = "m28"
$ws.Range("U16").Value = "m18"
$ws.Range("V16").Value = "m23/m18 Ratio:"

# Row 15 caption edited to describe the standardized output
$ws.Range("T15").Value = "And after standardizin we are getting:"

# New closing remark row
$ws.Range("T18").Value = "So the Mixed Reference Pattern coming out is as expected."

# --- Row 14 T:U become the "after standardizing" reference numbers (scientific style) ---
$ws.Range("T14").Value = 14.9849611
$ws.Range("U14").Value = 2.33229295
$ws.Range("T14:U14").NumberFormat = "0.00E+00"

# --- New row 17: the "after standardizing" measured ratio ---
$ws.Range("T17").Value = 23.56
$ws.Range("U17").Value = 3.667
$ws.Range("T17:U17").NumberFormat = "0.00E+00"
$ws.Range("V17").Formula = "=T17/U17"

# --- Fixed rounding problem in the tuning-correction fit coefficients (N and O columns) ---
$ws.Range("N2").Formula = "=(D2^2)*-0.00209891515351478+0.0415721558521753*D2+1.48645618255443"
$ws.Range("O2").Formula = "=0.00453957499583183*(D2^2) +D2*(-0.158786664417733) + 1.94715165380596"

$ws.Range("N3:N14").Formula = "=(D3^2)*-0.00209891515351478+0.0415721558521753*D3+1.48645618255443"
$ws.Range("O3:O14").Formula = "=0.00453957499583183*(D3^2) +D3*(-0.158786664417733) + 1.94715165380596"

# --- Selection moved to the new bottom-of-sheet content ---
$ws.Range("T19").Select()

$wb.Save()
